$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Stude" -> "Student Number" on the paragraph right after the
#    "<Project Name>" title.
# ------------------------------------------------------------------
$studeRange = $d.Content
$found = $studeRange.Find.Execute("Stude", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $studeRange.Collapse(0)
    $studeRange.InsertAfter("nt Number")
}

# ------------------------------------------------------------------
# 2) Add the "Problem Background" paragraphs right after the
#    "Background" Heading-2 paragraph (and before the following
#    blank paragraph that precedes "Scope").
# ------------------------------------------------------------------
$bgRange = $d.Content
$found = $bgRange.Find.Execute("Background", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $bgRange.Collapse(0)
    $bgRange.InsertAfter("`r<<<PARA1>>>`r<<<PARA2>>>")
}

$para1Text = "The number of automobiles on the road is increasing every day. As convenient as it makes the life of users, more cars running on the road simply means more chance of road accidents. One of the best ways to reduce the chances of accidents and loss of life and property in the future is to study when, how, and why the accidents have happened over the time at that area. We can thus see the patterns and trends of the accidents which helps a great deal in finding the leading causes of accidents and thus see where changes can be made."
$para2Text = "We have incorporated the use of data analysis and their tools to study injury and fatal crashes in Victoria based on various metrics."

$r1 = $d.Content
$r1.Find.Execute("<<<PARA1>>>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Style = "Normal"
$r1.ParagraphFormat.LeftIndent = 35.8
$r1.Text = $para1Text

$r2 = $d.Content
$r2.Find.Execute("<<<PARA2>>>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Style = "Normal"
$r2.ParagraphFormat.LeftIndent = 35.8
$r2.Text = $para2Text
